# Add 2022-Q3 data
# 1) Insert a new worksheet "2022-Q3" right after "总计" (pushing the existing
#    "2022-Q2" sheet to position 3).
# 2) Populate the new sheet with the Q3 per-fund holdings data.
# 3) Update the "总计" summary sheet: row 2 becomes the 2022-Q3 totals and a
#    new row 3 is added with the (previously row 2) 2022-Q2 totals.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: insert the new sheet. Worksheets.Add() with no args inserts the
# new sheet immediately after the active/first sheet, which is exactly the
# "总计","2022-Q3","2022-Q2" ordering we need.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add()
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Step 2: fill in the 2022-Q3 worksheet.
# Use Range.Copy(Destination) from cells on the summary sheet that already
# carry the desired cell style so the style gets reused (same style index)
# instead of a brand-new style being created.
# ---------------------------------------------------------------------

# Header row (B1:H1) - bold/centered style, same as summary sheet headers.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $summary.Cells.Item(1, 2).Copy($q3.Cells.Item(1, $col))
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Data rows 2-6.
$rows = @(
    @("009010", "华夏兴阳一年持有期混合", "26.58", "88.59", "2.82", "0.7496", 9),
    @("012584", "南方中国新兴经济9个月持有期混合（QDII）A", "2.69", "91.51", "2.91", "0.0783", 10),
    @("005255", "浦银安盛港股通量化混合A", "0.29", "78.68", "4.75", "0.0138", 6),
    @("012585", "南方中国新兴经济9个月持有期混合（QDII）C", "0.10", "91.51", "2.91", "0.0029", 10),
    @("013224", "浦银安盛港股通量化混合C", "0.05", "78.68", "4.75", "0.0024", 6)
)

$r = 2
foreach ($data in $rows) {
    # Column A: numeric index (0-based), styled like the summary's A2 (s=2).
    $summary.Cells.Item(2, 1).Copy($q3.Cells.Item($r, 1))
    $q3.Cells.Item($r, 1).Value = $r - 2

    # Column B: fund code, stored as text.
    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $data[0]
    $q3.Cells.Item($r, 2).Style = "Normal"

    # Column C: fund name (plain text, never looks numeric, fine as-is).
    $q3.Cells.Item($r, 3).Value = $data[1]

    # Columns D-G: numeric-looking figures, kept as text like the source file.
    for ($col = 4; $col -le 7; $col++) {
        $q3.Cells.Item($r, $col).NumberFormat = "@"
        $q3.Cells.Item($r, $col).Value = $data[$col - 2]
        $q3.Cells.Item($r, $col).Style = "Normal"
    }

    # Column H: rank, a real number.
    $q3.Cells.Item($r, 8).Value = $data[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 3: update the "总计" (summary) sheet.
# First, copy the formatting of A2 down into the new A3 row (so the
# "s=2" style is reused rather than a new style being minted), then fill
# in the values: row 3 gets what used to be in row 2 (2022-Q2 totals),
# and row 2 is overwritten with the new 2022-Q3 totals.
# ---------------------------------------------------------------------
$summary.Cells.Item(2, 1).Copy($summary.Cells.Item(3, 1))

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 2
$summary.Cells.Item(3, 4).Value = 0.03

$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 0.85
